# Task2Documentation.docx change log table update.
#
# The change log table's final row ("Added an analytics page...") was
# formatted slightly differently from every other row in the table: its
# first cell carried its w:cnfStyle marker inside the paragraph's w:pPr
# instead of the cell's w:tcPr (like all the other rows do), and its
# second cell had no w:cnfStyle/w:pPr at all. We fix that row to match the
# rest of the table's consistent formatting, and append a brand-new row
# underneath it documenting the dark/light mode rework, following that
# same consistent formatting.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the existing (malformed) last row and stash its text so we can
# restore it once it has been replaced by a correctly-formatted row.
$oldLastRow = $t.Rows.Item($t.Rows.Count)
$origDocText = $oldLastRow.Cells.Item(1).Range.Text
$origDateText = $oldLastRow.Cells.Item(2).Range.Text

# Rows.Add(beforeRow) clones the formatting of the row that sits just
# above the insertion point, which -- for every row except the malformed
# final one -- places w:cnfStyle on w:tcPr (col 1) / w:pPr (col 2) the
# way the rest of the table does. Insert two such well-formatted rows
# immediately above the malformed row, then delete the malformed row,
# leaving two clean rows in its place.
$refRow = $t.Rows.Item($t.Rows.Count)
$t.Rows.Add($refRow) | Out-Null

$refRow = $t.Rows.Item($t.Rows.Count)
$t.Rows.Add($refRow) | Out-Null

$t.Rows.Item($t.Rows.Count).Delete()

$fixedRow = $t.Rows.Item($t.Rows.Count - 1)
$newRow = $t.Rows.Item($t.Rows.Count)

# Restore the original changelog entry into the now-consistently-formatted
# row.
$fixedRow.Cells.Item(1).Range.Text = $origDocText
$fixedRow.Cells.Item(2).Range.Text = $origDateText

# Populate the new changelog entry.
$newRow.Cells.Item(1).Range.Text = "Changed how dark mode works, so that it’s easier to edit the light and dark mode colours and is more easily expandable"
$newRow.Cells.Item(2).Range.Text = "13/01/2025"
